$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 (was date 44574 block) -> becomes the 44223 block
$ws.Range("D4").Value2 = 44223
$ws.Range("M4").Value2 = 100
$ws.Range("N4").Value2 = 3500
$ws.Range("O4").Value2 = 4000
$ws.Range("P4").Value2 = 3750
$ws.Range("S4").Value2 = 1875

# Row 5 (was date 44574 block) -> becomes the 44223 block
$ws.Range("D5").Value2 = 44223
$ws.Range("M5").Value2 = 50
$ws.Range("N5").Value2 = 3000
$ws.Range("O5").Value2 = 3000
$ws.Range("P5").Value2 = 3000
$ws.Range("S5").Value2 = 1500

# Row 6 (was date 44223 block) -> becomes the 44574 block
$ws.Range("D6").Value2 = 44574
$ws.Range("M6").Value2 = 200
$ws.Range("N6").Value2 = 6000
$ws.Range("O6").Value2 = 7000
$ws.Range("P6").Value2 = 6500
$ws.Range("S6").Value2 = 3250

# Row 7 (was date 44223 block) -> becomes the 44574 block
$ws.Range("D7").Value2 = 44574
$ws.Range("M7").Value2 = 100
$ws.Range("N7").Value2 = 5000
$ws.Range("O7").Value2 = 5000
$ws.Range("P7").Value2 = 5000
$ws.Range("S7").Value2 = 2500
